# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '35.345.71'
$ws.Range('E2').Value2 = '  -0.56%  '
$ws.Range('D3').Value2 = '1.912.89'
$ws.Range('E3').Value2 = '  +0.21%  '
$ws.Range('E4').Value2 = '  -0.46%  '
$ws.Range('D5').Value2 = '''0.723'
$ws.Range('E5').Value2 = '  +10.29%  '
$ws.Range('D6').Value2 = '''252.54'
$ws.Range('E6').Value2 = '  +2.59%  '
$ws.Range('E7').Value2 = '  -0.28%  '
$ws.Range('D8').Value2 = '''40.52'
$ws.Range('E8').Value2 = '  -3.44%  '
$ws.Range('D9').Value2 = '''0.358'
$ws.Range('E9').Value2 = '  +3.39%  '
$ws.Range('D10').Value2 = '''52.74'
$ws.Range('E10').Value2 = '  +4.95%  '
$ws.Range('E11').Value2 = '  +1.98%  '
$ws.Range('D12').Value2 = '''0.0999'
$ws.Range('E12').Value2 = '  -0.30%  '
$ws.Range('D13').Value2 = '2.190.66'
$ws.Range('E13').Value2 = '  +0.14%  '
$ws.Range('D14').Value2 = '''12.55'
$ws.Range('E14').Value2 = '  +2.48%  '
$ws.Range('E15').Value2 = '  +1.98%  '
$ws.Range('D16').Value2 = '1.916.70'
$ws.Range('E16').Value2 = '  +0.67%  '
$ws.Range('D17').Value2 = '''4.88'
$ws.Range('E17').Value2 = '  -0.22%  '
$ws.Range('D18').Value2 = '35.341.97'
$ws.Range('E18').Value2 = '  -0.58%  '
$ws.Range('D19').Value2 = '''72.99'
$ws.Range('E19').Value2 = '  +0.87%  '
$ws.Range('E20').Value2 = '  +0.87%  '
$ws.Range('D21').Value2 = '''13.05'
$ws.Range('E21').Value2 = '  +3.29%  '
$ws.Range('D22').Value2 = '''241.55'
$ws.Range('E22').Value2 = '  -1.54%  '
$ws.Range('E23').Value2 = '  +5.02%  '
$ws.Range('E24').Value2 = '  -0.47%  '
$ws.Range('E25').Value2 = '  +1.25%  '
$ws.Range('D26').Value2 = '''2.33'
$ws.Range('E26').Value2 = '  +2.16%  '
$ws.Range('D27').Value2 = '''167.82'
$ws.Range('E27').Value2 = '  -1.91%  '
$ws.Range('D28').Value2 = '''8.67'
$ws.Range('E28').Value2 = '  +3.10%  '
$ws.Range('E29').Value2 = '  +5.41%  '
$ws.Range('E30').Value2 = '  +1.80%  '
$ws.Range('D31').Value2 = '4.132.19'
$ws.Range('E31').Value2 = '  +19.56%  '
$ws.Range('D32').Value2 = '''4.36'
$ws.Range('E32').Value2 = '  +4.43%  '
$ws.Range('E33').Value2 = '  +13.56%  '
$ws.Range('B34').Value2 = 'TrustWalletToken'
$ws.Range('C34').Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D34').Value2 = '''1.62'
$ws.Range('E34').Value2 = '  +20.48%  '
$ws.Range('B35').Value2 = 'Hedera'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value2 = '''0.0578'
$ws.Range('E35').Value2 = '  +1.41%  '
$ws.Range('D36').Value2 = '''4.24'
$ws.Range('E36').Value2 = '  +1.67%  '
$ws.Range('E37').Value2 = '  -0.39%  '
$ws.Range('D38').Value2 = '''0.912'
$ws.Range('E38').Value2 = '  -1.80%  '
$ws.Range('E39').Value2 = '  -1.41%  '
$ws.Range('D40').Value2 = '''17.47'
$ws.Range('E40').Value2 = '  +10.82%  '
$ws.Range('D41').Value2 = '''99.18'
$ws.Range('E41').Value2 = '  +8.76%  '
$ws.Range('E42').Value2 = '  +2.71%  '
$ws.Range('E43').Value2 = '  -1.45%  '
$ws.Range('D44').Value2 = '''0.0651'
$ws.Range('E44').Value2 = '  +2.58%  '
$ws.Range('D45').Value2 = '''2.48'
$ws.Range('E45').Value2 = '  +3.59%  '
$ws.Range('D46').Value2 = '1.348.09'
$ws.Range('E46').Value2 = '  -0.20%  '
$ws.Range('B47').Value2 = 'HuobiToken'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').Value2 = '''2.42'
$ws.Range('E47').Value2 = '  +0.48%  '
$ws.Range('B48').Value2 = 'FraxShare'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value2 = '''6.75'
$ws.Range('E48').Value2 = '  +2.78%  '
$ws.Range('E49').Value2 = '  -0.66%  '
$ws.Range('D50').Value2 = '''45.31'
$ws.Range('E50').Value2 = '  -4.63%  '
$ws.Range('D51').Value2 = '2.099.50'
$ws.Range('E51').Value2 = '  +0.09%  '
